# Updates the crypto price/volume table with refreshed coinranking.com data.
# Column D ("Price") holds decimal-looking numbers stored as TEXT in the
# original workbook (inline strings). Assigning a plain numeric-looking
# string to Range.Value lets Excel auto-convert it to a real number, so
# those assignments are prefixed with a literal leading apostrophe -- the
# normal Excel mechanism for forcing text entry -- to keep them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '34.764.57'

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.866.67'
$ws.Range('E3').Value = '  -2.52%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.92%  '

# Row 5 - BNB
$ws.Range('D5').Value = '''245.16'
$ws.Range('E5').Value = '  -3.47%  '

# Row 6 - XRP
$ws.Range('D6').Value = '''0.678'
$ws.Range('E6').Value = '  -6.03%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  -1.01%  '

# Row 8 - Solana
$ws.Range('D8').Value = '''41.65'
$ws.Range('E8').Value = '  +2.72%  '

# Row 9 - Cardano
$ws.Range('E9').Value = '  -3.93%  '

# Row 10 - Dogecoin
$ws.Range('D10').Value = '''0.0731'
$ws.Range('E10').Value = '  -2.78%  '

# Row 11 - TRON
$ws.Range('E11').Value = '  -2.56%  '

# Row 12 - Chainlink
$ws.Range('D12').Value = '''12.88'
$ws.Range('E12').Value = '  +1.30%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range('D13').Value = '2.138.29'
$ws.Range('E13').Value = '  -2.45%  '

# Row 14 - Polygon
$ws.Range('D14').Value = '''0.708'
$ws.Range('E14').Value = '  -1.49%  '

# Row 15 - WrappedEther
$ws.Range('D15').Value = '1.869.20'
$ws.Range('E15').Value = '  -2.49%  '

# Row 16 - Polkadot
$ws.Range('E16').Value = '  -2.08%  '

# Row 17 - WrappedBTC
$ws.Range('D17').Value = '34.745.25'
$ws.Range('E17').Value = '  -1.89%  '

# Row 18 - Litecoin
$ws.Range('D18').Value = '''72.12'
$ws.Range('E18').Value = '  -3.00%  '

# Row 19 - ShibaInu
$ws.Range('D19').Value = '0.0₃0809'
$ws.Range('E19').Value = '  -3.41%  '

# Row 20 - BitcoinCash
$ws.Range('D20').Value = '''242.40'
$ws.Range('E20').Value = '  -0.30%  '

# Row 21 - Avalanche
$ws.Range('D21').Value = '''12.55'
$ws.Range('E21').Value = '  -4.03%  '

# Row 22 - Uniswap
$ws.Range('E22').Value = '  -4.64%  '

# Row 23 - Dai
$ws.Range('E23').Value = '  -0.92%  '

# Row 24 - Toncoin
$ws.Range('D24').Value = '''2.47'
$ws.Range('E24').Value = '  +5.21%  '

# Row 25 - PancakeSwap
$ws.Range('D25').Value = '''2.15'
$ws.Range('E25').Value = '  -14.49%  '

# Row 26 - Monero
$ws.Range('D26').Value = '''163.21'
$ws.Range('E26').Value = '  -2.21%  '

# Row 27 - Cosmos
$ws.Range('D27').Value = '''8.33'
$ws.Range('E27').Value = '  -3.74%  '

# Row 28 - EthereumClassic
$ws.Range('D28').Value = '''18.05'
$ws.Range('E28').Value = '  -3.65%  '

# Row 29 - Stellar
$ws.Range('E29').Value = '  -5.59%  '

# Row 31 - TrustWalletToken
$ws.Range('D31').Value = '''1.71'
$ws.Range('E31').Value = '  +5.10%  '

# Row 32 - Filecoin
$ws.Range('E32').Value = '  -4.95%  '

# Row 33 - Hedera
$ws.Range('E33').Value = '  -1.31%  '

# Row 34 - BinanceUSD
$ws.Range('E34').Value = '  -0.96%  '

# Row 35 - InternetComputer(DFINITY)
$ws.Range('D35').Value = '''4.10'
$ws.Range('E35').Value = '  -2.88%  '

# Row 36 - ImmutableX
$ws.Range('D36').Value = '''0.827'
$ws.Range('E36').Value = '  -10.10%  '

# Row 37 - WEMIXToken
$ws.Range('D37').Value = '''1.58'
$ws.Range('E37').Value = '  -21.37%  '

# Row 38 - LidoDAOToken
$ws.Range('E38').Value = '  -3.06%  '

# Row 39 - Aave
$ws.Range('D39').Value = '''97.53'
$ws.Range('E39').Value = '  +0.18%  '

# Row 40 - Kaspa -> InjectiveProtocol
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').Value = '''16.94'
$ws.Range('E40').Value = '  -3.57%  '

# Row 41 - InjectiveProtocol -> Kaspa
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.0661'
$ws.Range('E41').Value = '  +1.08%  '

# Row 42 - VeChain
$ws.Range('D42').Value = '''0.0211'
$ws.Range('E42').Value = '  -1.15%  '

# Row 43 - ARBITRUM
$ws.Range('D43').Value = '''1.07'
$ws.Range('E43').Value = '  -4.75%  '

# Row 44 - Maker
$ws.Range('D44').Value = '1.281.15'
$ws.Range('E44').Value = '  -4.48%  '

# Row 45 - Cronos
$ws.Range('D45').Value = '''0.0816'
$ws.Range('E45').Value = '  +10.40%  '

# Row 46 - RenderToken
$ws.Range('D46').Value = '''2.29'
$ws.Range('E46').Value = '  -6.35%  '

# Row 48 - MXToken
$ws.Range('D48').Value = '''2.73'
$ws.Range('E48').Value = '  -1.46%  '

# Row 49 - Gas
$ws.Range('D49').Value = '''11.79'
$ws.Range('E49').Value = '  -0.71%  '

# Row 50 - FraxShare
$ws.Range('E50').Value = '  -7.08%  '

# Row 51 - MultiversX
$ws.Range('D51').Value = '''42.33'
$ws.Range('E51').Value = '  -5.32%  '
